$d = $word.ActiveDocument

# Update the date line (first paragraph) using Find/Replace - it is the only
# occurrence of this text in the document so a full-document search is safe.
$d.Content.Find.Execute("2024-11-15 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-16 Saturday", 2) | Out-Null

# Update each table cell by setting its Range.Text directly (row/column order).
# This avoids any cross-cell text collisions that a global Find/Replace could hit
# (several of the math expressions are substrings of other expressions in the grid).
$table = $d.Tables.Item(1)

$table.Cell(1, 1).Range.Text = "86-62="
$table.Cell(1, 2).Range.Text = "63+29="
$table.Cell(1, 3).Range.Text = "15+34="
$table.Cell(1, 4).Range.Text = "49+21="
$table.Cell(1, 5).Range.Text = "21+16="
$table.Cell(2, 1).Range.Text = "83-1="
$table.Cell(2, 2).Range.Text = "31+28="
$table.Cell(2, 3).Range.Text = "57+19="
$table.Cell(2, 4).Range.Text = "91+0="
$table.Cell(2, 5).Range.Text = "23+33="
$table.Cell(3, 1).Range.Text = "62-26="
$table.Cell(3, 2).Range.Text = "36-19="
$table.Cell(3, 3).Range.Text = "80+2="
$table.Cell(3, 4).Range.Text = "50-3="
$table.Cell(3, 5).Range.Text = "48+16="
$table.Cell(4, 1).Range.Text = "58+11="
$table.Cell(4, 2).Range.Text = "47+2="
$table.Cell(4, 3).Range.Text = "98-16="
$table.Cell(4, 4).Range.Text = "31+66="
$table.Cell(4, 5).Range.Text = "54+5="
$table.Cell(5, 1).Range.Text = "79+14="
$table.Cell(5, 2).Range.Text = "87-48="
$table.Cell(5, 3).Range.Text = "68-59="
$table.Cell(5, 4).Range.Text = "95-57="
$table.Cell(5, 5).Range.Text = "97-54="
$table.Cell(6, 1).Range.Text = "67-66="
$table.Cell(6, 2).Range.Text = "53-43="
$table.Cell(6, 3).Range.Text = "64+35="
$table.Cell(6, 4).Range.Text = "57-47="
$table.Cell(6, 5).Range.Text = "15-11="
$table.Cell(7, 1).Range.Text = "44+33="
$table.Cell(7, 2).Range.Text = "26+3="
$table.Cell(7, 3).Range.Text = "3+55="
$table.Cell(7, 4).Range.Text = "28-25="
$table.Cell(7, 5).Range.Text = "89+7="
$table.Cell(8, 1).Range.Text = "0+15="
$table.Cell(8, 2).Range.Text = "28-16="
$table.Cell(8, 3).Range.Text = "71-53="
$table.Cell(8, 4).Range.Text = "43+49="
$table.Cell(8, 5).Range.Text = "45+9="
$table.Cell(9, 1).Range.Text = "78-73="
$table.Cell(9, 2).Range.Text = "33-5="
$table.Cell(9, 3).Range.Text = "29+64="
$table.Cell(9, 4).Range.Text = "91-46="
$table.Cell(9, 5).Range.Text = "69+24="
$table.Cell(10, 1).Range.Text = "52-42="
$table.Cell(10, 2).Range.Text = "79+1="
$table.Cell(10, 3).Range.Text = "24+67="
$table.Cell(10, 4).Range.Text = "47+46="
$table.Cell(10, 5).Range.Text = "36+21="
$table.Cell(11, 1).Range.Text = "31+37="
$table.Cell(11, 2).Range.Text = "20+65="
$table.Cell(11, 3).Range.Text = "92-43="
$table.Cell(11, 4).Range.Text = "46+25="
$table.Cell(11, 5).Range.Text = "38+28="
$table.Cell(12, 1).Range.Text = "79-24="
$table.Cell(12, 2).Range.Text = "57-10="
$table.Cell(12, 3).Range.Text = "35+11="
$table.Cell(12, 4).Range.Text = "20+3="
$table.Cell(12, 5).Range.Text = "67+3="
$table.Cell(13, 1).Range.Text = "9+82="
$table.Cell(13, 2).Range.Text = "91-3="
$table.Cell(13, 3).Range.Text = "38+36="
$table.Cell(13, 4).Range.Text = "0+64="
$table.Cell(13, 5).Range.Text = "71+24="
$table.Cell(14, 1).Range.Text = "58-16="
$table.Cell(14, 2).Range.Text = "20-7="
$table.Cell(14, 3).Range.Text = "81+12="
$table.Cell(14, 4).Range.Text = "50+19="
$table.Cell(14, 5).Range.Text = "10+40="
$table.Cell(15, 1).Range.Text = "8+48="
$table.Cell(15, 2).Range.Text = "49-11="
$table.Cell(15, 3).Range.Text = "66-55="
$table.Cell(15, 4).Range.Text = "52+22="
$table.Cell(15, 5).Range.Text = "56-35="
$table.Cell(16, 1).Range.Text = "48+36="
$table.Cell(16, 2).Range.Text = "55-23="
$table.Cell(16, 3).Range.Text = "95-60="
$table.Cell(16, 4).Range.Text = "76-20="
$table.Cell(16, 5).Range.Text = "52-35="
$table.Cell(17, 1).Range.Text = "37+46="
$table.Cell(17, 2).Range.Text = "12+37="
$table.Cell(17, 3).Range.Text = "80-20="
$table.Cell(17, 4).Range.Text = "92-89="
$table.Cell(17, 5).Range.Text = "2+18="
$table.Cell(18, 1).Range.Text = "68-65="
$table.Cell(18, 2).Range.Text = "56-14="
$table.Cell(18, 3).Range.Text = "87+9="
$table.Cell(18, 4).Range.Text = "69+14="
$table.Cell(18, 5).Range.Text = "90-25="
$table.Cell(19, 1).Range.Text = "57-11="
$table.Cell(19, 2).Range.Text = "10-6="
$table.Cell(19, 3).Range.Text = "52+43="
$table.Cell(19, 4).Range.Text = "29+32="
$table.Cell(19, 5).Range.Text = "56-53="
$table.Cell(20, 1).Range.Text = "7-5="
$table.Cell(20, 2).Range.Text = "94-90="
$table.Cell(20, 3).Range.Text = "27+9="
$table.Cell(20, 4).Range.Text = "6+85="
$table.Cell(20, 5).Range.Text = "78-37="
